$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Full sentence text (bold lead-in "Meta description" + plain remainder).
$metaPara.Range.Text = "Meta description: Experience the epic online slot game Alexander’s Conquest. Play for free and read our review on the shifting game grid, graphics, bonuses, RTP, and more."

# Bold just the "Meta description" lead-in (splits off a second run for the
# remaining plain text).
$boldRange = $metaPara.Range.Duplicate
$boldRange.Find.ClearFormatting()
$boldRange.Find.Text = "Meta description"
$boldRange.Find.Execute() | Out-Null
$boldRange.Font.Bold = 1

# Leading empty run (matches the "<w:r/>, <w:r>text</w:r>" pattern used
# throughout the rest of the document).
$metaRange = $metaPara.Range.Duplicate
$metaRange.Collapse(1)
$metaRange.InsertBefore("")

# ---------------------------------------------------------------------------
# 2) Remove the trailing duplicate "Play Alexander's Conquest Free..." bold
#    paragraph near the end of the document (second-to-last paragraph).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupPara = $d.Paragraphs.Item($count - 1)
if ($dupPara.Range.Text -match "Play Alexander") {
    $dupPara.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 3) Replace the final (italic) paragraph's text with the new image prompt.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range.Duplicate
$lastRange.Find.ClearFormatting()
$lastRange.Find.Text = "Experience the epic online slot game Alexander’s Conquest. Play for free and read our review on the shifting game grid, graphics, bonuses, RTP, and more."
$lastRange.Find.Execute() | Out-Null
$lastRange.Text = "Prompt: Create a feature image for Alexander's Conquest that showcases a happy Maya warrior wearing glasses. The image should be in a cartoon style that suits the epic adventure theme of the game. Include the game's logo in the image and make sure the Maya warrior stands out as the main focal point. The background should feature a battlefield setting with soldiers and other relevant elements that reflect the game's storyline. Please use bold, vivid, and eye-catching colors to grab the viewer's attention and convey the excitement of playing the game."

Write-Output "done"
